$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New RAD test-run timestamps (Mon Oct 02 2023 run) replacing the previous
# "Fri Sep 29 ..." run's timestamps in column B, rows 2-29.
$dates = @(
    "Mon Oct 02 16:28:54 EDT 2023",
    "Mon Oct 02 16:29:04 EDT 2023",
    "Mon Oct 02 16:29:15 EDT 2023",
    "Mon Oct 02 16:29:29 EDT 2023",
    "Mon Oct 02 16:29:39 EDT 2023",
    "Mon Oct 02 16:29:49 EDT 2023",
    "Mon Oct 02 16:29:59 EDT 2023",
    "Mon Oct 02 16:30:09 EDT 2023",
    "Mon Oct 02 16:30:19 EDT 2023",
    "Mon Oct 02 16:30:30 EDT 2023",
    "Mon Oct 02 16:30:40 EDT 2023",
    "Mon Oct 02 16:30:50 EDT 2023",
    "Mon Oct 02 16:31:00 EDT 2023",
    "Mon Oct 02 16:31:10 EDT 2023",
    "Mon Oct 02 16:31:20 EDT 2023",
    "Mon Oct 02 16:31:31 EDT 2023",
    "Mon Oct 02 16:31:40 EDT 2023",
    "Mon Oct 02 16:31:51 EDT 2023",
    "Mon Oct 02 16:32:01 EDT 2023",
    "Mon Oct 02 16:32:11 EDT 2023",
    "Mon Oct 02 16:32:21 EDT 2023",
    "Mon Oct 02 16:32:31 EDT 2023",
    "Mon Oct 02 16:32:41 EDT 2023",
    "Mon Oct 02 16:32:51 EDT 2023",
    "Mon Oct 02 16:33:02 EDT 2023",
    "Mon Oct 02 16:33:12 EDT 2023",
    "Mon Oct 02 16:33:22 EDT 2023",
    "Mon Oct 02 16:33:32 EDT 2023"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $dates[$i]
}
